$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra car rows (4-7) that are no longer part of this report pull
$ws.Range("A4:O7").Delete()

# Update the remaining data row (row 3) to reflect the refreshed trace data
$ws.Range("C3").Value = "JOHNSTOWN"
$ws.Range("F3").Value = 19
$ws.Range("G3").Value = 1443
$ws.Range("H3").Value = "Placed Actual"
$ws.Range("I3").ClearContents()
$ws.Range("J3").Value = "JOHNSTOWN"

# Update the report description banner with the new pull date/time and event count
$ws.Range("A1").Value = "Description unknown, completed 10/24/2023 07:51:28 EDT, by WPJTOWN1.The search returned: 1 events."

# Fix up the selection so it only covers the remaining data cell
$null = $ws.Range("O3").Select()
